{"js": "// Update the multiplication problems in the practice-sheet table.\n// Each old problem string is unique in the document, so a simple\n// exact-match search/replace per pair reproduces the diff precisely\n// (only the <w:t> text of the matched run changes; all run/paragraph\n// formatting is left untouched).\nconst replacements = [\n  [\"90\u00d795=\", \"44\u00d791=\"],\n  [\"24\u00d722=\", \"44\u00d742=\"],\n  [\"72\u00d794=\", \"90\u00d781=\"],\n  [\"58\u00d714=\", \"86\u00d743=\"],\n  [\"62\u00d748=\", \"83\u00d730=\"],\n  [\"13\u00d720=\", \"67\u00d773=\"],\n  [\"77\u00d771=\", \"70\u00d793=\"],\n  [\"91\u00d745=\", \"69\u00d740=\"],\n  [\"97\u00d758=\", \"55\u00d776=\"],\n  [\"92\u00d799=\", \"16\u00d747=\"],\n  [\"51\u00d786=\", \"53\u00d723=\"],\n  [\"29\u00d759=\", \"98\u00d752=\"],\n  [\"14\u00d763=\", \"56\u00d790=\"],\n  [\"89\u00d758=\", \"71\u00d720=\"],\n  [\"32\u00d721=\", \"93\u00d723=\"],\n  [\"96\u00d749=\", \"89\u00d776=\"],\n  [\"62\u00d757=\", \"21\u00d745=\"],\n  [\"42\u00d731=\", \"92\u00d711=\"],\n  [\"19\u00d729=\", \"86\u00d764=\"],\n  [\"49\u00d781=\", \"74\u00d765=\"],\n  [\"25\u00d765=\", \"48\u00d726=\"],\n  [\"47\u00d712=\", \"73\u00d755=\"],\n  [\"43\u00d732=\", \"92\u00d768=\"],\n  [\"11\u00d760=\", \"89\u00d779=\"],\n  [\"79\u00d736=\", \"85\u00d780=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the multiplication problems in the practice-sheet table.\n# Each old problem string is unique in the document, so a plain\n# Find/Replace (MatchCase, no wildcards, ReplaceAll) per pair\n# reproduces the diff precisely - only the text of the matched run\n# changes, formatting is left untouched.\n\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$replacements = @(\n    @(\"90\u00d795=\", \"44\u00d791=\"),\n    @(\"24\u00d722=\", \"44\u00d742=\"),\n    @(\"72\u00d794=\", \"90\u00d781=\"),\n    @(\"58\u00d714=\", \"86\u00d743=\"),\n    @(\"62\u00d748=\", \"83\u00d730=\"),\n    @(\"13\u00d720=\", \"67\u00d773=\"),\n    @(\"77\u00d771=\", \"70\u00d793=\"),\n    @(\"91\u00d745=\", \"69\u00d740=\"),\n    @(\"97\u00d758=\", \"55\u00d776=\"),\n    @(\"92\u00d799=\", \"16\u00d747=\"),\n    @(\"51\u00d786=\", \"53\u00d723=\"),\n    @(\"29\u00d759=\", \"98\u00d752=\"),\n    @(\"14\u00d763=\", \"56\u00d790=\"),\n    @(\"89\u00d758=\", \"71\u00d720=\"),\n    @(\"32\u00d721=\", \"93\u00d723=\"),\n    @(\"96\u00d749=\", \"89\u00d776=\"),\n    @(\"62\u00d757=\", \"21\u00d745=\"),\n    @(\"42\u00d731=\", \"92\u00d711=\"),\n    @(\"19\u00d729=\", \"86\u00d764=\"),\n    @(\"49\u00d781=\", \"74\u00d765=\"),\n    @(\"25\u00d765=\", \"48\u00d726=\"),\n    @(\"47\u00d712=\", \"73\u00d755=\"),\n    @(\"43\u00d732=\", \"92\u00d768=\"),\n    @(\"11\u00d760=\", \"89\u00d779=\"),\n    @(\"79\u00d736=\", \"85\u00d780=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n}\n"}
